# Update the account-statement worksheet ("Estado de Cuenta") for
# NIT 9010159498: the previous data set (worker CRISPINA ISABEL DIAZ ARDILA
# with two overdue periods, plus worker HILDA ESTHER CONSUEGRA BELTRAN with
# one overdue period) is replaced with just the first "part" of the new
# statement: a single row for HILDA ESTHER CONSUEGRA BELTRAN / period 2408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table currently has 3 rows (16, 17, 18):
#   16: CC | 1047454117 | CRISPINA ISABEL DIAZ ARDILA | 2507 | 56940  | 1423500
#   17: CC | 1047454117 | CRISPINA ISABEL DIAZ ARDILA | 2506 | 56940  | 1423500
#   18: CC | 22391573   | HILDA ESTHER CONSUEGRA BELTRAN | 2408 | 16000 | 2000000
#
# Target: keep only the HILDA row, as the (only) data row, reusing the
# formatting that row 16 already has. So we overwrite row 16's values with
# row 18's values (format stays row 16's own), then delete rows 17 and 18.

$ws.Range("C16").Value = $ws.Range("C18").Value()
$ws.Range("D16").Value = $ws.Range("D18").Value()
$ws.Range("E16").Value = $ws.Range("E18").Value()
$ws.Range("F16").Value = $ws.Range("F18").Value()
$ws.Range("G16").Value = $ws.Range("G18").Value()

$ws.Rows("17:18").Delete()

# Update the summary KPIs to reflect the now-single remaining record.
$ws.Range("E11").Value = 16000   # VALOR MORA total
$ws.Range("C13").Value = 1       # Cant. Trabajadores
$ws.Range("F13").Value = 1       # Cant. Periodos
